$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header cells F1 and G1, matching the style of the existing header row (s="1")
$ws.Range("F1").Value = "conf.low"
$ws.Range("G1").Value = "conf.high"
$ws.Range("F1:G1").Font.Bold = $true
$ws.Range("F1:G1").HorizontalAlignment = $ws.Range("A1").HorizontalAlignment

# Fill in the new conf.low / conf.high numeric values for rows 2-12
$data = @{
    2  = @(0.8216757340518291, 0.9195105875976574)
    3  = @(-0.1088007211649279, -0.05550922784720759)
    4  = @(-0.1743482516071727, -0.1184970953753389)
    5  = @(-0.1021541559357472, -0.04654171088383176)
    6  = @(0.03952281396678342, 0.09466632125516136)
    7  = @(-0.02741228374835072, 0.0308315551021813)
    8  = @(-0.001107282487914878, 0.0005926691458439048)
    9  = @(-0.09379507011272285, -0.04742364612849419)
    10 = @(-0.4323415511637643, -0.3816422389321277)
    11 = @(-0.1790427958975642, -0.1278309946820979)
    12 = @(-0.01977429338400982, 0.03708617296711418)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Cells.Item($row, 6).Value = $vals[0]
    $ws.Cells.Item($row, 7).Value = $vals[1]
}

$wb.Save()
